# Weekly fruit/vegetable price update:
# - Insert 3 new rows at row 9 (shifts old rows 9-20 down to 12-23)
#   with a new week's worth of "Verde" Espárragos prices (fecha 44435).
# - Append 3 new rows (24-26) at the bottom with another new week's
#   worth of prices (fecha 44432).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 blank rows before row 9, pushing existing rows 9-20 down to 12-23 ---
$ws.Range("A9:R11").EntireRow.Insert()

# --- Step 2: populate the newly inserted rows 9-11 ---
# Common columns shared by the whole "Verde" Espárragos / Región Metropolitana block
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = 300000000
$commonG = "Espárragos"
$commonH = "Verde"
$commonO = "Región Metropolitana"
$commonR = "Hortaliza"

# Row 9: Banquete
$ws.Cells.Item(9, 1).Value = $commonA
$ws.Cells.Item(9, 2).Value = $commonB
$ws.Cells.Item(9, 3).Value = $commonC
$ws.Cells.Item(9, 4).Value = 44435
$ws.Cells.Item(9, 5).Value = $commonE
$ws.Cells.Item(9, 6).Value = $commonF
$ws.Cells.Item(9, 7).Value = $commonG
$ws.Cells.Item(9, 8).Value = $commonH
$ws.Cells.Item(9, 9).Value = "Banquete"
$ws.Cells.Item(9, 10).Value = 7
$ws.Cells.Item(9, 11).Value = 38000
$ws.Cells.Item(9, 12).Value = 40000
$ws.Cells.Item(9, 13).Value = 39143
$ws.Cells.Item(9, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 15).Value = $commonO
$ws.Cells.Item(9, 16).Value = 3914
$ws.Cells.Item(9, 17).Value = 10
$ws.Cells.Item(9, 18).Value = $commonR

# Row 10: Primera
$ws.Cells.Item(10, 1).Value = $commonA
$ws.Cells.Item(10, 2).Value = $commonB
$ws.Cells.Item(10, 3).Value = $commonC
$ws.Cells.Item(10, 4).Value = 44435
$ws.Cells.Item(10, 5).Value = $commonE
$ws.Cells.Item(10, 6).Value = $commonF
$ws.Cells.Item(10, 7).Value = $commonG
$ws.Cells.Item(10, 8).Value = $commonH
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 34000
$ws.Cells.Item(10, 12).Value = 36000
$ws.Cells.Item(10, 13).Value = 34960
$ws.Cells.Item(10, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 15).Value = $commonO
$ws.Cells.Item(10, 16).Value = 3496
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(10, 18).Value = $commonR

# Row 11: Segunda
$ws.Cells.Item(11, 1).Value = $commonA
$ws.Cells.Item(11, 2).Value = $commonB
$ws.Cells.Item(11, 3).Value = $commonC
$ws.Cells.Item(11, 4).Value = 44435
$ws.Cells.Item(11, 5).Value = $commonE
$ws.Cells.Item(11, 6).Value = $commonF
$ws.Cells.Item(11, 7).Value = $commonG
$ws.Cells.Item(11, 8).Value = $commonH
$ws.Cells.Item(11, 9).Value = "Segunda"
$ws.Cells.Item(11, 10).Value = 16
$ws.Cells.Item(11, 11).Value = 30000
$ws.Cells.Item(11, 12).Value = 32000
$ws.Cells.Item(11, 13).Value = 31000
$ws.Cells.Item(11, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(11, 15).Value = $commonO
$ws.Cells.Item(11, 16).Value = 3100
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(11, 18).Value = $commonR

# --- Step 3: append 3 new rows (24-26) at the bottom ---

# Row 24: Banquete
$ws.Cells.Item(24, 1).Value = $commonA
$ws.Cells.Item(24, 2).Value = $commonB
$ws.Cells.Item(24, 3).Value = $commonC
$ws.Cells.Item(24, 4).Value = 44432
$ws.Cells.Item(24, 5).Value = $commonE
$ws.Cells.Item(24, 6).Value = $commonF
$ws.Cells.Item(24, 7).Value = $commonG
$ws.Cells.Item(24, 8).Value = $commonH
$ws.Cells.Item(24, 9).Value = "Banquete"
$ws.Cells.Item(24, 10).Value = 7
$ws.Cells.Item(24, 11).Value = 38000
$ws.Cells.Item(24, 12).Value = 40000
$ws.Cells.Item(24, 13).Value = 39143
$ws.Cells.Item(24, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(24, 15).Value = $commonO
$ws.Cells.Item(24, 16).Value = 3914
$ws.Cells.Item(24, 17).Value = 10
$ws.Cells.Item(24, 18).Value = $commonR

# Row 25: Primera
$ws.Cells.Item(25, 1).Value = $commonA
$ws.Cells.Item(25, 2).Value = $commonB
$ws.Cells.Item(25, 3).Value = $commonC
$ws.Cells.Item(25, 4).Value = 44432
$ws.Cells.Item(25, 5).Value = $commonE
$ws.Cells.Item(25, 6).Value = $commonF
$ws.Cells.Item(25, 7).Value = $commonG
$ws.Cells.Item(25, 8).Value = $commonH
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 25
$ws.Cells.Item(25, 11).Value = 34000
$ws.Cells.Item(25, 12).Value = 36000
$ws.Cells.Item(25, 13).Value = 34960
$ws.Cells.Item(25, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(25, 15).Value = $commonO
$ws.Cells.Item(25, 16).Value = 3496
$ws.Cells.Item(25, 17).Value = 10
$ws.Cells.Item(25, 18).Value = $commonR

# Row 26: Segunda
$ws.Cells.Item(26, 1).Value = $commonA
$ws.Cells.Item(26, 2).Value = $commonB
$ws.Cells.Item(26, 3).Value = $commonC
$ws.Cells.Item(26, 4).Value = 44432
$ws.Cells.Item(26, 5).Value = $commonE
$ws.Cells.Item(26, 6).Value = $commonF
$ws.Cells.Item(26, 7).Value = $commonG
$ws.Cells.Item(26, 8).Value = $commonH
$ws.Cells.Item(26, 9).Value = "Segunda"
$ws.Cells.Item(26, 10).Value = 16
$ws.Cells.Item(26, 11).Value = 30000
$ws.Cells.Item(26, 12).Value = 32000
$ws.Cells.Item(26, 13).Value = 31000
$ws.Cells.Item(26, 14).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(26, 15).Value = $commonO
$ws.Cells.Item(26, 16).Value = 3100
$ws.Cells.Item(26, 17).Value = 10
$ws.Cells.Item(26, 18).Value = $commonR

# Ensure column D keeps the date/datetime number format used elsewhere in the column
$ws.Range("D9:D11").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("D24:D26").NumberFormat = $ws.Range("D8").NumberFormat
